$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zufriedenheit")

# Update header row
$ws.Range("A1").Value = "vpn"
$ws.Range("B1").Value = "messzeitpunkt"
$ws.Range("C1").Value = "zufriedenheit"

# Update VP identifiers to lowercase (column A, rows 2-9)
$ws.Range("A2").Value = "vp_1"
$ws.Range("A3").Value = "vp_2"
$ws.Range("A4").Value = "vp_3"
$ws.Range("A5").Value = "vp_4"
$ws.Range("A6").Value = "vp_1"
$ws.Range("A7").Value = "vp_2"
$ws.Range("A8").Value = "vp_3"
$ws.Range("A9").Value = "vp_4"

# Set active cell to A2 as in the final file's selection
$ws.Range("A2").Select()
